$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-11, columns B-G (ME, MAE, MSE, RMSE, SE, N)
# A naive forecaster bug fix shifted the error-table window down by one
# quarter: row 2 now holds freshly computed figures, and rows 3-11 take
# on the values that previously lived one row above them.
$data = @(
    @(0.1783908196033299, 0.3606156554386025, 0.2599511937740667, 0.5098540906711122, 0.4943913024279584, 15),
    @(0.2147957804815247, 0.3580819848962541, 0.2088379301833887, 0.4569878884427777, 0.4185882721160311, 14),
    @(0.245473867636748,  0.3683843928281823, 0.2358514199634523, 0.4856453644002507, 0.4361500127557235, 13),
    @(0.301463903713677,  0.3567050099425655, 0.2583510523008795, 0.5082824532687308, 0.4274285484895888, 12),
    @(0.2900334324511264, 0.3439303588265401, 0.1874043702602189, 0.4329022640968963, 0.3370659818980856, 11),
    @(0.2833661826906564, 0.2909807410604855, 0.2099604875429734, 0.4582144558424291, 0.3795671424286398, 10),
    @(0.3070978867771534, 0.3733636115588901, 0.2272551859511885, 0.476712896774556,  0.3867354821117653, 9),
    @(0.3125739411747067, 0.4101416856701894, 0.3161064439746084, 0.5622334425971194, 0.499604386659364,  8),
    @(0.2627214152773011, 0.3235434349789264, 0.1669924024756296, 0.4086470389904099, 0.3380801243818421, 7),
    @(0.3023013051968123, 0.3856149542229345, 0.3827274330185637, 0.618649685216572,  0.5912779588939882, 6)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Count; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}
